# "Updated for week 8"
#
# Adds the two new Week 8 timesheet entries (2/25/2019 and 2/27/2019),
# which pushes the Week 8 weekly total from 0 to 4.75 hours and the
# cumulative project total (which cascades through Weeks 9-15 and the
# Final summary sheet via formulas) from 41.83 to 46.58. It also moves
# the "active"/selected worksheet from Week 7 to Week 8.

$wb  = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("Week 7")
$ws8 = $wb.Worksheets.Item("Week 8")

# --- New Week 8 data rows -------------------------------------------------

# Row 2: Mon 2/25/2019, 4:00 PM - 6:00 PM, "Created remaining pages", 2 hrs
$ws8.Range("A2").Value2 = 43521
$ws8.Range("A2").NumberFormat = "m/d/yy"
$ws8.Range("B2").Value2 = 0.66666666666666663
$ws8.Range("B2").NumberFormat = "h:mm AM/PM"
$ws8.Range("C2").Value2 = 0.75
$ws8.Range("C2").NumberFormat = "h:mm AM/PM"
$ws8.Range("D2").Value2 = "Created remaining pages"
$ws8.Range("E2").Value2 = 2
$ws8.Rows("2:2").RowHeight = 18

# Row 3: Wed 2/27/2019, 12:30 PM - 3:15 PM, "Worked on page content and CSS", 2.75 hrs
$ws8.Range("A3").Value2 = 43523
$ws8.Range("A3").NumberFormat = "m/d/yy"
$ws8.Range("B3").Value2 = 0.52083333333333337
$ws8.Range("B3").NumberFormat = "h:mm AM/PM"
$ws8.Range("C3").Value2 = 0.63541666666666663
$ws8.Range("C3").NumberFormat = "h:mm AM/PM"
$ws8.Range("D3").Value2 = "Worked on page content and CSS"
$ws8.Range("E3").Value2 = 2.75
$ws8.Rows("3:3").RowHeight = 18

# --- Active sheet / selection moves from Week 7 to Week 8 ----------------

$ws7.Activate()
$ws7.Range("C6").Select()

$ws8.Activate()
$ws8.Range("D4").Select()
